# Auto-generated Excel COM-interop script to apply scheduled-runner updates
# to the Adamantoise_Profits workbook (per-job profit recalculation sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 1634.4762
$ws.Range("I28").Value = 1804.1111
$ws.Range("J28").Value = 616.6667
$ws.Range("K28").Value = 1804.1111
$ws.Range("L28").Value = 616.6667
$ws.Range("M28").Value = -1319.1111
$ws.Range("N28").Value = -1586.6667
# Row 101
$ws.Range("H101").Value = 519.2
$ws.Range("I101").Value = 568.3333
$ws.Range("K101").Value = 1704.9999
$ws.Range("M101").Value = -82.99990000000003
# Row 119
$ws.Range("H119").Value = 3750
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 3750
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 11250
$ws.Range("M119").ClearContents()
$ws.Range("N119").Value = -20926
# Row 132
$ws.Range("H132").Value = 4003.8223
$ws.Range("I132").Value = 4198.108
$ws.Range("K132").Value = 12594.324
$ws.Range("M132").Value = -10064.324
# Row 138
$ws.Range("H138").Value = 2779.2432
$ws.Range("I138").Value = 1978.2
$ws.Range("J138").Value = 3187.9387
$ws.Range("K138").Value = 5934.6
$ws.Range("L138").Value = 9563.8161
$ws.Range("M138").Value = -794.6000000000004
$ws.Range("N138").Value = -19843.8161

$ws = $wb.Worksheets.Item("ARM")
# Row 44
$ws.Range("H44").Value = 66970
$ws.Range("J44").Value = 66970
$ws.Range("L44").Value = 66970
$ws.Range("N44").Value = -67946
# Row 61
$ws.Range("H61").Value = 2501.6785
$ws.Range("I61").Value = 2217.5217
$ws.Range("J61").Value = 3808.8
$ws.Range("K61").Value = 2217.5217
$ws.Range("L61").Value = 3808.8
$ws.Range("M61").Value = -2005.5217
$ws.Range("N61").Value = -4232.8
# Row 63
$ws.Range("H63").Value = 3659.3
$ws.Range("I63").Value = 2199.75
$ws.Range("J63").Value = 4632.3335
$ws.Range("K63").Value = 2199.75
$ws.Range("L63").Value = 4632.3335
$ws.Range("M63").Value = -1513.75
$ws.Range("N63").Value = -6004.3335
# Row 66
$ws.Range("H66").Value = 3659.3
$ws.Range("I66").Value = 2199.75
$ws.Range("J66").Value = 4632.3335
$ws.Range("K66").Value = 10998.75
$ws.Range("L66").Value = 23161.6675
$ws.Range("M66").Value = -7566.75
$ws.Range("N66").Value = -30025.6675
# Row 74
$ws.Range("H74").Value = 1998.8
$ws.Range("I74").Value = 1623.5
$ws.Range("J74").Value = 3500
$ws.Range("K74").Value = 1623.5
$ws.Range("L74").Value = 3500
$ws.Range("M74").Value = -749.5
$ws.Range("N74").Value = -5248
# Row 77
$ws.Range("H77").Value = 1998.8
$ws.Range("I77").Value = 1623.5
$ws.Range("J77").Value = 3500
$ws.Range("K77").Value = 8117.5
$ws.Range("L77").Value = 17500
$ws.Range("M77").Value = -3749.5
$ws.Range("N77").Value = -26236
# Row 80
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
# Row 82
$ws.Range("H82").Value = 60000
$ws.Range("J82").Value = 60000
$ws.Range("L82").Value = 60000
$ws.Range("N82").Value = -60722
# Row 83
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
# Row 85
$ws.Range("H85").Value = 60000
$ws.Range("J85").Value = 60000
$ws.Range("L85").Value = 60000
$ws.Range("N85").Value = -62496
# Row 136
$ws.Range("H136").Value = 2501.6785
$ws.Range("I136").Value = 2217.5217
$ws.Range("J136").Value = 3808.8
$ws.Range("K136").Value = 6652.5651
$ws.Range("L136").Value = 11426.4
$ws.Range("M136").Value = -4102.5651
$ws.Range("N136").Value = -16526.4

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 36802944
$ws.Range("I134").Value = 10207294
$ws.Range("J134").Value = 83345336
$ws.Range("K134").Value = 30621882
$ws.Range("L134").Value = 250036008
$ws.Range("M134").Value = -30619347
$ws.Range("N134").Value = -250041078

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3013.3547
$ws.Range("I31").Value = 2163.9048
$ws.Range("J31").Value = 4797.2
$ws.Range("K31").Value = 2163.9048
$ws.Range("L31").Value = 4797.2
$ws.Range("M31").Value = -1868.9048
$ws.Range("N31").Value = -5387.2
# Row 34
$ws.Range("H34").Value = 3013.3547
$ws.Range("I34").Value = 2163.9048
$ws.Range("J34").Value = 4797.2
$ws.Range("K34").Value = 2163.9048
$ws.Range("L34").Value = 4797.2
$ws.Range("M34").Value = -1961.9048
$ws.Range("N34").Value = -5201.2
# Row 52
$ws.Range("H52").Value = 20104.5
$ws.Range("J52").Value = 29500
$ws.Range("L52").Value = 29500
$ws.Range("N52").Value = -30088
# Row 132
$ws.Range("H132").Value = 1510.8
$ws.Range("I132").Value = 1119.9412
$ws.Range("J132").Value = 2341.375
$ws.Range("K132").Value = 3359.8236
$ws.Range("L132").Value = 7024.125
$ws.Range("M132").Value = -829.8235999999997
$ws.Range("N132").Value = -12084.125
# Row 134
$ws.Range("H134").Value = 2768.8333
$ws.Range("I134").Value = 1854.4
$ws.Range("J134").Value = 3911.875
$ws.Range("K134").Value = 5563.200000000001
$ws.Range("L134").Value = 11735.625
$ws.Range("M134").Value = -3028.200000000001
$ws.Range("N134").Value = -16805.625

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1869.36
$ws.Range("I102").Value = 1547.05
$ws.Range("J102").Value = 3158.6
$ws.Range("K102").Value = 1547.05
$ws.Range("L102").Value = 3158.6
$ws.Range("M102").Value = 74.95000000000005
$ws.Range("N102").Value = -6402.6
# Row 132
$ws.Range("H132").Value = 1645.3077
$ws.Range("I132").Value = 1682.4166
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 5047.2498
$ws.Range("L132").Value = 3600
$ws.Range("M132").Value = -2517.2498
$ws.Range("N132").Value = -8660

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 3504.4666
$ws.Range("I132").Value = 3132
$ws.Range("K132").Value = 9396
$ws.Range("M132").Value = -6866
# Row 136
$ws.Range("H136").Value = 4052
$ws.Range("I136").Value = 3802
$ws.Range("J136").Value = 4352
$ws.Range("K136").Value = 11406
$ws.Range("L136").Value = 13056
$ws.Range("M136").Value = -8856
$ws.Range("N136").Value = -18156

$ws = $wb.Worksheets.Item("WVR")
# Row 47
$ws.Range("H47").Value = 39999
$ws.Range("I47").Value = 39999
$ws.Range("K47").Value = 39999
$ws.Range("M47").Value = -39427
# Row 132
$ws.Range("H132").Value = 3419.5334
$ws.Range("I132").Value = 3108.2856
$ws.Range("J132").Value = 7777
$ws.Range("K132").Value = 9324.856800000001
$ws.Range("L132").Value = 23331
$ws.Range("M132").Value = -6794.856800000001
$ws.Range("N132").Value = -28391
# Row 136
$ws.Range("H136").Value = 3231.9285
$ws.Range("I136").Value = 1409.6666
$ws.Range("K136").Value = 4228.9998
$ws.Range("M136").Value = -1678.9998
